# 17/07/2018 MAMATHA CHICK IN
#
# 1) Merge the two runs that make up the last existing timestamp
#    ("SUN Jul 15" + " 13:56:33 IST 2018") into a single run, by doing a
#    Find/Replace of the full visible string with itself (Word coalesces
#    the run when it rewrites the match).
# 2) Insert a brand-new purchase-record block (dated "MON Jul 16 13:15:10
#    IST 2018") right after the last record's "Amount balance" paragraph,
#    matching the exact paragraph/run layout used throughout the log.

$d = $word.ActiveDocument

# --- Step 1: coalesce "SUN Jul 15" / " 13:56:33 IST 2018" into one run ---
[void]$d.Content.Find.Execute("SUN Jul 15 13:56:33 IST 2018", $false, $false, $false, `
    $false, $false, $true, 1, $false, "SUN Jul 15 13:56:33 IST 2018", 2)

# --- Step 2: locate the end of the last record ("Amount balance - 41604.0") ---
$rng = $d.Content
[void]$rng.Find.Execute("41604.0", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lastPara = $rng.Paragraphs(1)

# Insert point sits just before that paragraph's own paragraph mark, so the
# new content becomes its own paragraph(s) rather than merging into the
# (already-existing) empty paragraph that follows it.
$insertAt = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)

$rPr = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr>'
$rPrB = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/></w:rPr>'
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newBlock = @"
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPrB</w:pPr></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPrB</w:pPr></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPrB</w:pPr></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPr</w:pPr><w:r>$rPr<w:t>MON Jul 16</w:t></w:r><w:r>$rPr<w:t xml:space="preserve"> 13:15:10 IST 2018</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPr</w:pPr><w:r>$rPr<w:t>Person Name</w:t></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/><w:t>- MB</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPr</w:pPr><w:r>$rPr<w:t>Bill number</w:t></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/><w:t>- 2631</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPr</w:pPr><w:r>$rPr<w:t>---------------------------------------------------------------</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPr</w:pPr><w:r>$rPr<w:t>Item Name</w:t></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/><w:t>- CARROT</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPr</w:pPr><w:r>$rPr<w:t>Number of Pockets</w:t></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/><w:t>- 2</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPr</w:pPr><w:r>$rPr<w:t>Number of KGs</w:t></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/><w:t>- 173</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPr</w:pPr><w:r>$rPr<w:t>Rate</w:t></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/><w:t>- 26</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPr</w:pPr><w:r>$rPr<w:t>Transport &amp; Miscellaneous</w:t></w:r><w:r>$rPr<w:tab/><w:t>- 40</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPr</w:pPr><w:r>$rPr<w:t>Total Price</w:t></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/></w:r><w:r>$rPr<w:tab/><w:t>- 4538.0</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPrB</w:pPr><w:r>$rPrB<w:t>Amount balance</w:t></w:r><w:r>$rPrB<w:tab/></w:r><w:r>$rPrB<w:tab/></w:r><w:r>$rPrB<w:tab/><w:t>- 46142.0</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPr</w:pPr></w:p>
<w:p $ns><w:pPr><w:pStyle w:val="PlainText"/>$rPrB</w:pPr></w:p>
"@

$insertAt.InsertXML($newBlock)
